$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "433.81") are stored as literal text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '66.601.53'
$ws.Cells.Item(2, 5).Value = '  -0.96%  '
$ws.Cells.Item(3, 4).Value = '3.790.71'
$ws.Cells.Item(3, 5).Value = '  -1.35%  '
$ws.Cells.Item(4, 5).Value = '  -0.45%  '
$ws.Cells.Item(5, 4).Value = '433.81'
$ws.Cells.Item(5, 5).Value = '  +4.95%  '
$ws.Cells.Item(6, 4).Value = '139.88'
$ws.Cells.Item(6, 5).Value = '  +4.80%  '
$ws.Cells.Item(7, 5).Value = '  +0.44%  '
$ws.Cells.Item(8, 4).Value = '0.999'
$ws.Cells.Item(8, 5).Value = '  -0.22%  '
$ws.Cells.Item(9, 4).Value = '0.736'
$ws.Cells.Item(9, 5).Value = '  -1.52%  '
$ws.Cells.Item(10, 5).Value = '  -12.16%  '
$ws.Cells.Item(11, 5).Value = '  -16.86%  '
$ws.Cells.Item(12, 4).Value = '42.63'
$ws.Cells.Item(12, 5).Value = '  +3.05%  '
$ws.Cells.Item(13, 5).Value = '  +3.70%  '
$ws.Cells.Item(14, 4).Value = '4.419.60'
$ws.Cells.Item(14, 5).Value = '  -1.14%  '
$ws.Cells.Item(15, 4).Value = '15.09'
$ws.Cells.Item(15, 5).Value = '  +1.87%  '
$ws.Cells.Item(16, 5).Value = '  -0.33%  '
$ws.Cells.Item(17, 4).Value = '3.803.48'
$ws.Cells.Item(17, 5).Value = '  -0.08%  '
$ws.Cells.Item(18, 5).Value = '  +1.55%  '
$ws.Cells.Item(19, 5).Value = '  +3.54%  '
$ws.Cells.Item(20, 4).Value = '66.672.58'
$ws.Cells.Item(20, 5).Value = '  -1.48%  '
$ws.Cells.Item(21, 4).Value = '411.30'
$ws.Cells.Item(21, 5).Value = '  -1.86%  '
$ws.Cells.Item(22, 4).Value = '14.67'
$ws.Cells.Item(22, 5).Value = '  -1.87%  '
$ws.Cells.Item(23, 5).Value = '  +5.66%  '
$ws.Cells.Item(24, 4).Value = '85.30'
$ws.Cells.Item(24, 5).Value = '  -1.61%  '
$ws.Cells.Item(25, 4).Value = '36.89'
$ws.Cells.Item(25, 5).Value = '  +0.34%  '
$ws.Cells.Item(26, 5).Value = '  +5.05%  '
$ws.Cells.Item(27, 2).Value = 'LEO'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(27, 4).Value = '5.61'
$ws.Cells.Item(27, 5).Value = '  -2.11%  '
$ws.Cells.Item(28, 2).Value = 'RenderToken'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(28, 4).Value = '9.70'
$ws.Cells.Item(28, 5).Value = '  +33.60%  '
$ws.Cells.Item(29, 4).Value = '9.80'
$ws.Cells.Item(29, 5).Value = '  +2.12%  '
$ws.Cells.Item(30, 5).Value = '  +11.36%  '
$ws.Cells.Item(31, 2).Value = 'Bittensor'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(31, 4).Value = '717.65'
$ws.Cells.Item(31, 5).Value = '  +2.65%  '
$ws.Cells.Item(32, 2).Value = 'Cosmos'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(32, 4).Value = '13.89'
$ws.Cells.Item(32, 5).Value = '  +10.73%  '
$ws.Cells.Item(33, 5).Value = '  +0.70%  '
$ws.Cells.Item(34, 4).Value = '41.52'
$ws.Cells.Item(34, 5).Value = '  +5.72%  '
$ws.Cells.Item(35, 4).Value = '0.999'
$ws.Cells.Item(35, 5).Value = '  -0.02%  '
$ws.Cells.Item(36, 4).Value = '5.66'
$ws.Cells.Item(36, 5).Value = '  +27.96%  '
$ws.Cells.Item(37, 4).Value = '0.149'
$ws.Cells.Item(37, 5).Value = '  -3.63%  '
$ws.Cells.Item(38, 4).Value = '55.77'
$ws.Cells.Item(38, 5).Value = '  +0.21%  '
$ws.Cells.Item(39, 4).Value = '0.0473'
$ws.Cells.Item(39, 5).Value = '  +1.79%  '
$ws.Cells.Item(40, 4).Value = '2.75'
$ws.Cells.Item(40, 5).Value = '  +38.72%  '
$ws.Cells.Item(41, 4).Value = '2.96'
$ws.Cells.Item(41, 5).Value = '  -4.44%  '
$ws.Cells.Item(42, 4).Value = '0.0₃0696'
$ws.Cells.Item(42, 5).Value = '  -12.52%  '
$ws.Cells.Item(43, 5).Value = '  +3.30%  '
$ws.Cells.Item(44, 5).Value = '  +0.10%  '
$ws.Cells.Item(45, 4).Value = '3.23'
$ws.Cells.Item(45, 5).Value = '  +2.91%  '
$ws.Cells.Item(46, 5).Value = '  +8.40%  '
$ws.Cells.Item(47, 5).Value = '  +0.02%  '
$ws.Cells.Item(48, 5).Value = '  +3.14%  '
$ws.Cells.Item(49, 4).Value = '2.09'
$ws.Cells.Item(49, 5).Value = '  -1.25%  '
$ws.Cells.Item(50, 4).Value = '142.15'
$ws.Cells.Item(50, 5).Value = '  -4.07%  '
$ws.Cells.Item(51, 4).Value = '2.82'
$ws.Cells.Item(51, 5).Value = '  -1.70%  '

# Restore the default (unstyled) cell style now that the text is committed,
# so formatting matches the original workbook.
$ws.Range("D2:D51").Style = "Normal"
